$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D sometimes hold plain decimal-looking text (e.g. '619.71').
# Force those specific cells to Text format first so Excel doesn't coerce them
# into numeric cells (and lose trailing zeros / introduce float noise).
$ws.Range('D2').Value = '97.557.35'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '3.298.50'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '254.47'
$ws.Range('E5').Value = '  +4.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '619.71'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.43'
$ws.Range('E7').Value = '  +26.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.397'
$ws.Range('E8').Value = '  +2.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  +13.93%  '
$ws.Range('D11').Value = '3.297.35'
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.52'
$ws.Range('E13').Value = '  +10.04%  '
$ws.Range('D14').Value = '97.340.77'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000245'
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').Value = '3.920.49'
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = '3.302.62'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.50'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.09'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.17'
$ws.Range('E21').Value = '  +5.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '475.97'
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.36'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000203'
$ws.Range('E24').Value = '  -2.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.56'
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.89'
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.78'
$ws.Range('E27').Value = '  -2.51%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.297'
$ws.Range('E28').Value = '  +22.51%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.474.95'
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  +2.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.131'
$ws.Range('E32').Value = '  +8.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.76'
$ws.Range('E33').Value = '  +5.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.41'
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.147'
$ws.Range('E36').Value = '  -2.88%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.12'
$ws.Range('E37').Value = '  -3.71%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '487.84'
$ws.Range('E40').Value = '  -2.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.452'
$ws.Range('E41').Value = '  +0.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.68'
$ws.Range('E42').Value = '  +6.75%  '
$ws.Range('E43').Value = '  -3.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.793'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.13'
$ws.Range('E46').Value = '  -4.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '157.87'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.837'
$ws.Range('E49').Value = '  +4.64%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.11'
$ws.Range('E50').Value = '  +11.96%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.60'
$ws.Range('E51').Value = '  +1.17%  '
